{"js": "// \"Restored the level 6 depth header layout\"\n//\n// 1. Drop the stray `_GoBack` bookmark that sat on the title paragraph\n//    (Word relocates `_GoBack` to the spot of the last edit, which\n//    renumbers every other bookmark id down by one).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Re-create `_GoBack` around \"uur \" inside \"... in figuur ...\", which is\n//    where the edit actually happened \u2014 this splits that run into three.\nconst goBackTarget = context.document.body.search(\"uur \", { matchCase: true });\ngoBackTarget.load(\"items\");\nawait context.sync();\n\nif (goBackTarget.items.length > 0) {\n  goBackTarget.items[0].insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 3. Restore the old \"Heading 6\" look used for the level-6 depth header:\n//    no forced page break, dark-blue \"Text 2\" theme color, 11pt text.\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal\");\nawait context.sync();\n\nconst heading6 = styles.items.find((s) => s.nameLocal === \"Heading 6\");\nif (heading6) {\n  heading6.paragraphFormat.pageBreakBefore = false;\n  heading6.font.color = \"#1F497D\";\n  heading6.font.size = 11;\n  await context.sync();\n}\n", "ps1": "# \"Restored the level 6 depth header layout\"\n\n$d = $word.ActiveDocument\n\n# 1. Drop the stray `_GoBack` bookmark that sat on the title paragraph\n#    (Word relocates `_GoBack` to the spot of the last edit, which\n#    renumbers every other bookmark id down by one).\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# 2. Re-create `_GoBack` around \"uur \" inside \"... in figuur ...\", which is\n#    where the edit actually happened \u2014 this splits that run into three.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"uur \"\n$found = $find.Execute()\nif ($found) {\n    $d.Bookmarks.Add(\"_GoBack\", $rng)\n}\n\n# 3. Restore the old \"Heading 6\" look used for the level-6 depth header:\n#    no forced page break, dark-blue \"Text 2\" theme color, 11pt text.\n$st = $d.Styles(\"Heading 6\")\n$st.ParagraphFormat.PageBreakBefore = $false\n$st.Font.TextColor.ObjectThemeColor = 15\n$st.Font.Size = 11\n"}
